$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - update "想去人数" (F) column values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 292
$ws1.Range("F4").Value = 2411
$ws1.Range("F5").Value = 1780
$ws1.Range("F6").Value = 341
$ws1.Range("F8").Value = 833
$ws1.Range("F9").Value = 167

# Sheet "全部类型" (sheet4) - update "想去人数" (F) column values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 292
$ws4.Range("F4").Value = 2411
$ws4.Range("F5").Value = 1780
$ws4.Range("F6").Value = 341
$ws4.Range("F9").Value = 833
$ws4.Range("F10").Value = 167
